# Apply the 24-May-2024 cryptos list refresh (prices & 1h volume %, plus two row swaps).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.967.50'
$ws.Range('E2').Value = '  +2.50%  '
# Row 3
$ws.Range('D3').Value = '3.744.97'
$ws.Range('E3').Value = '  +0.82%  '
# Row 4
$ws.Range('E4').Value = '  +0.41%  '
# Row 5
$ws.Range('D5').Value = "'601.66"
$ws.Range('E5').Value = '  +2.12%  '
# Row 6
$ws.Range('D6').Value = "'167.75"
$ws.Range('E6').Value = '  -1.52%  '
# Row 7
$ws.Range('D7').Value = '3.743.10'
$ws.Range('E7').Value = '  +0.72%  '
# Row 8
$ws.Range('E8').Value = '  +0.02%  '
# Row 9
$ws.Range('E9').Value = '  +3.58%  '
# Row 10
$ws.Range('D10').Value = "'0.165"
$ws.Range('E10').Value = '  +5.58%  '
# Row 11
$ws.Range('E11').Value = '  +3.24%  '
# Row 12
$ws.Range('D12').Value = "'0.460"
$ws.Range('E12').Value = '  +0.63%  '
# Row 13
$ws.Range('D13').Value = "'38.20"
$ws.Range('E13').Value = '  +2.67%  '
# Row 14
$ws.Range('D14').Value = "'0.0000244"
$ws.Range('E14').Value = '  +2.10%  '
# Row 15
$ws.Range('D15').Value = '4.372.24'
$ws.Range('E15').Value = '  +0.81%  '
# Row 16
$ws.Range('D16').Value = '3.752.75'
$ws.Range('E16').Value = '  +1.48%  '
# Row 17
$ws.Range('D17').Value = '68.974.96'
$ws.Range('E17').Value = '  +2.73%  '
# Row 18
$ws.Range('D18').Value = "'7.27"
$ws.Range('E18').Value = '  +2.66%  '
# Row 19
$ws.Range('E19').Value = '  +1.15%  '
# Row 20
$ws.Range('D20').Value = "'17.26"
$ws.Range('E20').Value = '  +8.09%  '
# Row 21
$ws.Range('D21').Value = "'497.59"
$ws.Range('E21').Value = '  +3.06%  '
# Row 22
$ws.Range('D22').Value = "'10.25"
$ws.Range('E22').Value = '  +14.89%  '
# Row 23
$ws.Range('D23').Value = "'0.724"
$ws.Range('E23').Value = '  +2.41%  '
# Row 24
$ws.Range('D24').Value = "'85.26"
$ws.Range('E24').Value = '  +2.88%  '
# Row 25
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Value = "'2.30"
$ws.Range('E25').Value = '  -0.73%  '
# Row 26
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = "'0.0000141"
$ws.Range('E26').Value = '  +3.16%  '
# Row 27
$ws.Range('D27').Value = "'12.29"
$ws.Range('E27').Value = '  +1.63%  '
# Row 28
$ws.Range('E28').Value = '  +0.03%  '
# Row 29
$ws.Range('E29').Value = '  +0.22%  '
# Row 30
$ws.Range('E30').Value = '  +2.06%  '
# Row 31
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = "'2.42"
$ws.Range('E31').Value = '  +2.13%  '
# Row 32
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = "'7.95"
$ws.Range('E32').Value = '  +4.27%  '
# Row 33
$ws.Range('D33').Value = "'31.76"
$ws.Range('E33').Value = '  -0.92%  '
# Row 34
$ws.Range('D34').Value = '3.887.21'
$ws.Range('E34').Value = '  +0.95%  '
# Row 35
$ws.Range('E35').Value = '  +1.47%  '
# Row 36
$ws.Range('D36').Value = '3.675.88'
$ws.Range('E36').Value = '  +0.65%  '
# Row 37
$ws.Range('E37').Value = '  +0.71%  '
# Row 38
$ws.Range('E38').Value = '  +2.52%  '
# Row 39
$ws.Range('D39').Value = "'5.80"
$ws.Range('E39').Value = '  +2.84%  '
# Row 40
$ws.Range('D40').Value = "'0.134"
$ws.Range('E40').Value = '  +0.34%  '
# Row 41
$ws.Range('E41').Value = '  +1.37%  '
# Row 42
$ws.Range('D42').Value = "'439.44"
$ws.Range('E42').Value = '  -0.98%  '
# Row 43
$ws.Range('D43').Value = "'49.07"
$ws.Range('E43').Value = '  +0.98%  '
# Row 44
$ws.Range('D44').Value = "'1.99"
$ws.Range('E44').Value = '  +1.01%  '
# Row 45
$ws.Range('D45').Value = "'2.87"
$ws.Range('E45').Value = '  +2.76%  '
# Row 46
$ws.Range('E46').Value = '  +2.53%  '
# Row 47
$ws.Range('E47').Value = '  -0.01%  '
# Row 48
$ws.Range('D48').Value = "'40.44"
$ws.Range('E48').Value = '  -1.75%  '
# Row 49
$ws.Range('D49').Value = "'142.64"
$ws.Range('E49').Value = '  +1.53%  '
# Row 50
$ws.Range('E50').Value = '  +2.26%  '
# Row 51
$ws.Range('D51').Value = '2.749.79'
$ws.Range('E51').Value = '  -0.94%  '
